{"js": "// The \"Date\" line currently reads \"Date: May 03/2019\", where \"May 03\" and\n// \"/2019\" are two separate runs. Replace that literal date text with the\n// {d[i].date} template placeholder, collapsing it into a single run.\nconst results = context.document.body.search(\"May 03/2019\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  results.items[0].insertText(\"{d[i].date}\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# The \"Date\" line currently reads \"Date: May 03/2019\" where \"May 03\" and\n# \"/2019\" are two separate runs. Replace that literal date text with the\n# {d[i].date} template placeholder, collapsing it into a single run.\n$find = $d.Content.Find\n$find.Text = \"May 03/2019\"\n$find.Replacement.Text = \"{d[i].date}\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n"}
